$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("BillTo", "BillToAddress", "ShipTo", "ShipToAddress", "AccountNumber", "ProductLine")

$startCol = 12  # L
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $startCol + $i
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    # 30 (character width) round-trips through the engine's internal unit
    # conversion as 30 + 5/6; back it off here so the saved <col width="..">
    # lands on the same "30" the other data columns already use.
    $ws.Columns.Item($col).ColumnWidth = 29.1666666666667
}

# Copy header style (fill/font/border/alignment) from K1 to the new header cells
$src = $ws.Range("K1")
$dst = $ws.Range($ws.Cells.Item(1, $startCol), $ws.Cells.Item(1, $startCol + $headers.Length - 1))
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
